$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 10529210
$ws.Range("I116").Value = 20002566
$ws.Range("J116").Value = 3258.889
$ws.Range("K116").Value = 20002566
$ws.Range("L116").Value = 3258.889
$ws.Range("M116").Value = -19999124
$ws.Range("N116").Value = -10142.889
$ws.Range("H125").Value = 72240.42999999999
$ws.Range("J125").Value = 960.1429000000001
$ws.Range("L125").Value = 8641.286100000001
$ws.Range("N125").Value = -13561.2861
$ws.Range("H129").Value = 1027.2059
$ws.Range("J129").Value = 1217.48
$ws.Range("L129").Value = 3652.44
$ws.Range("N129").Value = -13652.44
$ws.Range("H137").Value = 1890.258
$ws.Range("I137").Value = 1769.5
$ws.Range("J137").Value = 2109.818
$ws.Range("K137").Value = 5308.5
$ws.Range("L137").Value = 6329.454000000001
$ws.Range("M137").Value = -2758.5
$ws.Range("N137").Value = -11429.454

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12917.955
$ws.Range("I32").Value = 13254.9
$ws.Range("J32").Value = 10029.857
$ws.Range("K32").Value = 13254.9
$ws.Range("L32").Value = 10029.857
$ws.Range("M32").Value = -12967.9
$ws.Range("N32").Value = -10603.857
$ws.Range("H45").Value = 1673
$ws.Range("I45").Value = 1673
$ws.Range("K45").Value = 1673
$ws.Range("M45").Value = -1296
$ws.Range("H61").Value = 1442.0667
$ws.Range("I61").Value = 1312.6666
$ws.Range("J61").Value = 1959.6666
$ws.Range("K61").Value = 1312.6666
$ws.Range("L61").Value = 1959.6666
$ws.Range("M61").Value = -1100.6666
$ws.Range("N61").Value = -2383.6666
$ws.Range("H63").Value = 2619.25
$ws.Range("J63").Value = 5395
$ws.Range("L63").Value = 5395
$ws.Range("N63").Value = -6767
$ws.Range("H66").Value = 2619.25
$ws.Range("J66").Value = 5395
$ws.Range("L66").Value = 26975
$ws.Range("N66").Value = -33839
$ws.Range("H74").Value = 1081.8125
$ws.Range("I74").Value = 878.0909
$ws.Range("K74").Value = 878.0909
$ws.Range("M74").Value = -4.090900000000033
$ws.Range("H77").Value = 1081.8125
$ws.Range("I77").Value = 878.0909
$ws.Range("K77").Value = 4390.4545
$ws.Range("M77").Value = -22.45449999999983
$ws.Range("H122").Value = 1898.96
$ws.Range("I122").Value = 2000.381
$ws.Range("K122").Value = 6001.143
$ws.Range("M122").Value = -3551.143
$ws.Range("H136").Value = 1442.0667
$ws.Range("I136").Value = 1312.6666
$ws.Range("J136").Value = 1959.6666
$ws.Range("K136").Value = 3937.9998
$ws.Range("L136").Value = 5878.9998
$ws.Range("M136").Value = -1387.9998
$ws.Range("N136").Value = -10978.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 916.6667
$ws.Range("I29").Value = 916.6667
$ws.Range("K29").Value = 916.6667
$ws.Range("M29").Value = -627.6667
$ws.Range("H99").Value = 2110
$ws.Range("I99").Value = 1166.6666
$ws.Range("J99").Value = 2276.4707
$ws.Range("K99").Value = 1166.6666
$ws.Range("L99").Value = 2276.4707
$ws.Range("M99").Value = 331.3334
$ws.Range("N99").Value = -5272.4707
$ws.Range("H134").Value = 2045.4546
$ws.Range("I134").Value = 1870.9143
$ws.Range("J134").Value = 2724.2222
$ws.Range("K134").Value = 5612.742899999999
$ws.Range("L134").Value = 8172.6666
$ws.Range("M134").Value = -3077.742899999999
$ws.Range("N134").Value = -13242.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23258830
$ws.Range("I31").Value = 38463870
$ws.Range("J31").Value = 4057.7646
$ws.Range("K31").Value = 38463870
$ws.Range("L31").Value = 4057.7646
$ws.Range("M31").Value = -38463575
$ws.Range("N31").Value = -4647.7646
$ws.Range("H34").Value = 23258830
$ws.Range("I34").Value = 38463870
$ws.Range("J34").Value = 4057.7646
$ws.Range("K34").Value = 38463870
$ws.Range("L34").Value = 4057.7646
$ws.Range("M34").Value = -38463668
$ws.Range("N34").Value = -4461.7646
$ws.Range("H58").Value = 1249.9736
$ws.Range("I58").Value = 1212.1212
$ws.Range("J58").Value = 1499.8
$ws.Range("K58").Value = 1212.1212
$ws.Range("L58").Value = 1499.8
$ws.Range("M58").Value = -1009.1212
$ws.Range("N58").Value = -1905.8
$ws.Range("H132").Value = 1952.1708
$ws.Range("I132").Value = 1688.9429
$ws.Range("J132").Value = 3487.6667
$ws.Range("K132").Value = 5066.8287
$ws.Range("L132").Value = 10463.0001
$ws.Range("M132").Value = -2536.8287
$ws.Range("N132").Value = -15523.0001
$ws.Range("H134").Value = 1543.0344
$ws.Range("I134").Value = 1398.2
$ws.Range("J134").Value = 1864.8889
$ws.Range("K134").Value = 4194.6
$ws.Range("L134").Value = 5594.6667
$ws.Range("M134").Value = -1659.6
$ws.Range("N134").Value = -10664.6667
$ws.Range("H136").Value = 1249.9736
$ws.Range("I136").Value = 1212.1212
$ws.Range("J136").Value = 1499.8
$ws.Range("K136").Value = 3636.3636
$ws.Range("L136").Value = 4499.4
$ws.Range("M136").Value = -1086.3636
$ws.Range("N136").Value = -9599.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1632.0968
$ws.Range("I5").Value = 2731.923
$ws.Range("J5").Value = 837.7778
$ws.Range("K5").Value = 8195.769
$ws.Range("L5").Value = 2513.3334
$ws.Range("M5").Value = -8083.769
$ws.Range("N5").Value = -2737.3334
$ws.Range("H129").Value = 1352533.6
$ws.Range("I129").Value = 475.36365
$ws.Range("J129").Value = 1924558.2
$ws.Range("K129").Value = 1426.09095
$ws.Range("L129").Value = 5773674.6
$ws.Range("M129").Value = 3573.90905
$ws.Range("N129").Value = -5783674.6
$ws.Range("H135").Value = 1632.0968
$ws.Range("I135").Value = 2731.923
$ws.Range("J135").Value = 837.7778
$ws.Range("K135").Value = 24587.307
$ws.Range("L135").Value = 7540.000199999999
$ws.Range("M135").Value = -22052.307
$ws.Range("N135").Value = -12610.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 19990
$ws.Range("I53").Value = 10000
$ws.Range("K53").Value = 10000
$ws.Range("M53").Value = -9369
$ws.Range("H113").Value = 1567.9445
$ws.Range("I113").Value = 1724.3334
$ws.Range("J113").Value = 1255.1666
$ws.Range("K113").Value = 1724.3334
$ws.Range("L113").Value = 1255.1666
$ws.Range("M113").Value = 445.6666
$ws.Range("N113").Value = -5595.1666
$ws.Range("H122").Value = 2852.9167
$ws.Range("I122").Value = 3241.9583
$ws.Range("K122").Value = 9725.874899999999
$ws.Range("M122").Value = -7275.874899999999
$ws.Range("H123").Value = 18748.182
$ws.Range("J123").Value = 18748.182
$ws.Range("L123").Value = 18748.182
$ws.Range("N123").Value = -23648.182
$ws.Range("H132").Value = 1813.3096
$ws.Range("I132").Value = 1436.5151
$ws.Range("J132").Value = 3194.889
$ws.Range("K132").Value = 4309.5453
$ws.Range("L132").Value = 9584.667000000001
$ws.Range("M132").Value = -1779.5453
$ws.Range("N132").Value = -14644.667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 13000
$ws.Range("J44").Value = 13000
$ws.Range("L44").Value = 13000
$ws.Range("M44").Value = -13912
$ws.Range("H69").Value = 500163
$ws.Range("J69").Value = 500163
$ws.Range("L69").Value = 500163
$ws.Range("N69").Value = -501785
$ws.Range("H72").Value = 500163
$ws.Range("J72").Value = 500163
$ws.Range("L72").Value = 1500489
$ws.Range("N72").Value = -1508601
$ws.Range("H132").Value = 5512.364
$ws.Range("I132").Value = 5540.5186
$ws.Range("K132").Value = 16621.5558
$ws.Range("M132").Value = -14091.5558

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 94332.63
$ws.Range("I81").Value = 203699.8
$ws.Range("J81").Value = 3193.3333
$ws.Range("K81").Value = 407399.6
$ws.Range("L81").Value = 6386.6666
$ws.Range("M81").Value = -406338.6
$ws.Range("N81").Value = -8508.6666
$ws.Range("H84").Value = 94332.63
$ws.Range("I84").Value = 203699.8
$ws.Range("J84").Value = 3193.3333
$ws.Range("K84").Value = 2036998
$ws.Range("L84").Value = 31933.333
$ws.Range("M84").Value = -2031694
$ws.Range("N84").Value = -42541.333
$ws.Range("H125").Value = 61264.23
$ws.Range("J125").Value = 61264.23
$ws.Range("L125").Value = 61264.23
$ws.Range("N125").Value = -71104.23000000001
$ws.Range("H132").Value = 1870.6444
$ws.Range("I132").Value = 1526.88
$ws.Range("J132").Value = 2300.35
$ws.Range("K132").Value = 4580.64
$ws.Range("L132").Value = 6901.049999999999
$ws.Range("M132").Value = -2050.64
$ws.Range("N132").Value = -11961.05
